{"js": "// Replace every arithmetic-problem cell in the table with its updated\n// problem text. The table has 20 rows x 5 columns = 100 cells; every\n// cell's text is replaced according to its position in the table\n// (matching the order the problems appear in the document), since some\n// problem strings (e.g. \"28+25=\") occur more than once before the edit\n// but map to different replacements afterward.\nconst newRows = [\n  [\"17+18=\", \"70-56=\", \"73-18=\", \"74-67=\", \"39+18=\"],\n  [\"54-7=\", \"8+47=\", \"82-37=\", \"91-47=\", \"60-22=\"],\n  [\"83-68=\", \"33+19=\", \"95-58=\", \"55+8=\", \"39+26=\"],\n  [\"61-22=\", \"77+14=\", \"64-17=\", \"38+58=\", \"85-76=\"],\n  [\"50-32=\", \"28+48=\", \"22+69=\", \"5+28=\", \"73-27=\"],\n  [\"95-69=\", \"75-49=\", \"75+16=\", \"56+39=\", \"27+54=\"],\n  [\"9+74=\", \"19+28=\", \"72-59=\", \"26+59=\", \"72-16=\"],\n  [\"51-9=\", \"90-21=\", \"16+46=\", \"84-68=\", \"63+28=\"],\n  [\"70-5=\", \"8+36=\", \"95-39=\", \"8+83=\", \"9+47=\"],\n  [\"55-16=\", \"7+54=\", \"85-49=\", \"65-26=\", \"52-26=\"],\n  [\"19+66=\", \"74-39=\", \"49+46=\", \"62-46=\", \"66-7=\"],\n  [\"18+49=\", \"96-57=\", \"76-68=\", \"32-25=\", \"26-17=\"],\n  [\"95-56=\", \"14+28=\", \"57+6=\", \"16+15=\", \"17+26=\"],\n  [\"57-28=\", \"91-3=\", \"87+9=\", \"74-48=\", \"50-21=\"],\n  [\"24-19=\", \"93-35=\", \"74-15=\", \"46-8=\", \"65+6=\"],\n  [\"73-4=\", \"26+25=\", \"36+16=\", \"77-49=\", \"7+9=\"],\n  [\"17+75=\", \"6+65=\", \"42-35=\", \"23+38=\", \"90-12=\"],\n  [\"91-33=\", \"84-46=\", \"94-87=\", \"37+24=\", \"52-35=\"],\n  [\"80-58=\", \"36+27=\", \"15+9=\", \"90-47=\", \"27-8=\"],\n  [\"36+58=\", \"25+49=\", \"90-77=\", \"39+54=\", \"47-18=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst currentRowCount = table.values.length;\nconst currentColCount = table.values.length ? table.values[0].length : 0;\nif (currentRowCount !== newRows.length || currentColCount !== newRows[0].length) {\n  throw new Error(\n    \"Unexpected table shape: \" + currentRowCount + \"x\" + currentColCount\n  );\n}\n\ntable.values = newRows;\nawait context.sync();\n", "ps1": "# Replace every arithmetic-problem cell in the table with its updated\n# problem text. The table has 20 rows x 5 columns = 100 cells; every\n# cell's text is replaced according to its position in the table\n# (matching the order the problems appear in the document), since some\n# problem strings (e.g. \"28+25=\") occur more than once before the edit\n# but map to different replacements afterward.\n$newRows = @(\n  @(\"17+18=\", \"70-56=\", \"73-18=\", \"74-67=\", \"39+18=\"),\n  @(\"54-7=\", \"8+47=\", \"82-37=\", \"91-47=\", \"60-22=\"),\n  @(\"83-68=\", \"33+19=\", \"95-58=\", \"55+8=\", \"39+26=\"),\n  @(\"61-22=\", \"77+14=\", \"64-17=\", \"38+58=\", \"85-76=\"),\n  @(\"50-32=\", \"28+48=\", \"22+69=\", \"5+28=\", \"73-27=\"),\n  @(\"95-69=\", \"75-49=\", \"75+16=\", \"56+39=\", \"27+54=\"),\n  @(\"9+74=\", \"19+28=\", \"72-59=\", \"26+59=\", \"72-16=\"),\n  @(\"51-9=\", \"90-21=\", \"16+46=\", \"84-68=\", \"63+28=\"),\n  @(\"70-5=\", \"8+36=\", \"95-39=\", \"8+83=\", \"9+47=\"),\n  @(\"55-16=\", \"7+54=\", \"85-49=\", \"65-26=\", \"52-26=\"),\n  @(\"19+66=\", \"74-39=\", \"49+46=\", \"62-46=\", \"66-7=\"),\n  @(\"18+49=\", \"96-57=\", \"76-68=\", \"32-25=\", \"26-17=\"),\n  @(\"95-56=\", \"14+28=\", \"57+6=\", \"16+15=\", \"17+26=\"),\n  @(\"57-28=\", \"91-3=\", \"87+9=\", \"74-48=\", \"50-21=\"),\n  @(\"24-19=\", \"93-35=\", \"74-15=\", \"46-8=\", \"65+6=\"),\n  @(\"73-4=\", \"26+25=\", \"36+16=\", \"77-49=\", \"7+9=\"),\n  @(\"17+75=\", \"6+65=\", \"42-35=\", \"23+38=\", \"90-12=\"),\n  @(\"91-33=\", \"84-46=\", \"94-87=\", \"37+24=\", \"52-35=\"),\n  @(\"80-58=\", \"36+27=\", \"15+9=\", \"90-47=\", \"27-8=\"),\n  @(\"36+58=\", \"25+49=\", \"90-77=\", \"39+54=\", \"47-18=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nif ($t.Rows.Count -ne $newRows.Count -or $t.Columns.Count -ne $newRows[0].Count) {\n    throw (\"Unexpected table shape: \" + $t.Rows.Count + \"x\" + $t.Columns.Count)\n}\n\nfor ($r = 0; $r -lt $newRows.Count; $r++) {\n    $row = $newRows[$r]\n    for ($c = 0; $c -lt $row.Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n    }\n}\n"}
